$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOB1258")

# Semestre ideal: EA-6 -> EA-7
$ws.Range("B9").Value = "EA-7"
$ws.Range("C9").Value = "EA-7"

# Requisitos: replace first requirement text, remove second requirement row
$ws.Range("B25").Value = "LOB1217 -  Operações Unitárias e Processos  (Requisito fraco)`n"
$ws.Range("C25").Value = "LOB1217 -  Operações Unitárias e Processos  (Requisito fraco)`n"

# Remove the now-obsolete last requirement row entirely
$ws.Rows("26:26").Delete()
